# The sheet originally has an extra leading column A (row-index-ish values
# that duplicate column F), a header row in B1:F1, and data in A2:F5.
# The edit removes that leading column A entirely (shifting B:F -> A:E),
# and renames the "MODEL_CONDITION" header text to "MODELCONDITION".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire first column; everything to the right shifts left.
$ws.Columns.Item(1).Delete()

# Fix the header text (was "MODEL_CONDITION", now "MODELCONDITION").
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")
